# Update the "Background" slide's last bullet paragraph (Content Placeholder 2)
# from the two-run "I am in this course to learn best practices for establishing
# reproducible data workflow pipelines" to a single new sentence about the
# student's eDNA metabarcoding pipeline project.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item("Content Placeholder 2")
$tr = $sh.TextFrame.TextRange

$oldText = "I am in this course to learn best practices for establishing reproducible data workflow pipelines"
$newText = "I am in the beginning stages of establishing an eDNA metabarcoding pipeline, and am in this class to learn how to create a manageable and reproducible project structure"

$fullText = $tr.Text
$idx = $fullText.IndexOf($oldText)

if ($idx -lt 0) {
    throw "Could not locate target paragraph text in shape."
}

# TextRange.Characters is 1-based; $idx from .IndexOf is 0-based, and since
# .Text already includes the paragraph-mark characters, the mapping is direct.
$startPos = $idx + 1
$len = $oldText.Length

$target = $tr.Characters($startPos, $len)
$target.Text = $newText
